$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.458.24'
$ws.Range('E2').Value = '  +1.07%  '

$ws.Range('D3').Value = '2.329.64'
$ws.Range('E3').Value = '  +1.07%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.08'
$ws.Range('E5').Value = '  +1.01%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.37'
$ws.Range('E6').Value = '  +0.65%  '

$ws.Range('E7').Value = '  -0.74%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.505'
$ws.Range('E9').Value = '  -0.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.82'
$ws.Range('E10').Value = '  -0.18%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.48'
$ws.Range('E11').Value = '  +7.20%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0800'
$ws.Range('E12').Value = '  +1.25%  '

$ws.Range('E13').Value = '  +0.69%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.93'
$ws.Range('E14').Value = '  +1.96%  '

$ws.Range('D15').Value = '2.694.99'
$ws.Range('E15').Value = '  +1.17%  '

$ws.Range('D16').Value = '2.316.08'
$ws.Range('E16').Value = '  +0.54%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').Value = '  +1.51%  '

$ws.Range('D18').Value = '43.364.80'
$ws.Range('E18').Value = '  +1.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.80'
$ws.Range('E19').Value = '  +1.39%  '

$ws.Range('E20').Value = '  -0.11%  '

$ws.Range('E21').Value = '  +1.03%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.12'
$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.77'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.27'
$ws.Range('E24').Value = '  +6.09%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  +0.37%  '

$ws.Range('E26').Value = '  +0.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.14'
$ws.Range('E27').Value = '  -1.22%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.07'
$ws.Range('E28').Value = '  +0.75%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '164.60'
$ws.Range('E29').Value = '  -0.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.15'
$ws.Range('E30').Value = '  +1.02%  '

$ws.Range('E31').Value = '  +0.56%  '

$ws.Range('E32').Value = '  -0.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.02'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.91'
$ws.Range('E34').Value = '  +5.53%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.50'
$ws.Range('E35').Value = '  -6.97%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0707'
$ws.Range('E36').Value = '  +2.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.35'
$ws.Range('E37').Value = '  -1.31%  '

$ws.Range('E38').Value = '  -0.15%  '

$ws.Range('E39').Value = '  +1.10%  '

$ws.Range('E40').Value = '  +2.14%  '

$ws.Range('E41').Value = '  -0.10%  '

$ws.Range('D42').Value = '1.990.73'
$ws.Range('E42').Value = '  -0.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.69'
$ws.Range('E43').Value = '  +6.23%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.88'
$ws.Range('E44').Value = '  +7.93%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0282'
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('E46').Value = '  -1.45%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.80'
$ws.Range('E47').Value = '  +0.29%  '

$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.36'
$ws.Range('E48').Value = '  +1.40%  '

$ws.Range('D49').Value = '2.561.11'
$ws.Range('E49').Value = '  +1.19%  '

$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.87'
$ws.Range('E50').Value = '  -3.59%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.03'
$ws.Range('E51').Value = '  +1.06%  '

